# Added grant addTag to * privilege for Case and Complaint.
#
# The "Access Control Rules" table on Sheet1 gets two new rows inserted:
#   - A new COMPLAINT rule "Complaint – Anybody can add tag" /
#     "grant addTag to *" inserted right after the existing row for
#     "Complaint – Anybody can subscribe" (i.e. just before the first
#     CASE_FILE rule, which pushes it - and everything below - down by
#     one row).
#   - A new CASE_FILE rule "Case File – anyone can add tag" /
#     "grant addTag to *" inserted right after the existing row for
#     "Case File – anyone can subscribe" (i.e. just before the first
#     TASK rule).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New COMPLAINT row: "Complaint – Anybody can add tag" ---------------
# Insert a blank row at 30 (old row 30 "Case File – Assignee Read Access"
# and everything below shifts down by one).
$null = $ws.Rows.Item(30).Insert()

# Row 29 ("Complaint – Anybody can subscribe") has exactly the formatting
# we want for the new row (COMPLAINT rule, same column styles), so copy
# its formats down into the freshly inserted row 30.
$null = $ws.Range("A29:G29").Copy()
$null = $ws.Range("A30:G30").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("B30").Value = "Complaint – Anybody can add tag"
$ws.Range("C30").Value = "COMPLAINT"
$ws.Range("G30").Value = "grant addTag to *"
$ws.Rows.Item(30).RowHeight = 30

# --- New CASE_FILE row: "Case File – anyone can add tag" ----------------
# After the first insertion, the old row 44 ("Case File – anyone can
# subscribe") now lives at row 45. Insert a new blank row right after it,
# at row 46 (shifting the TASK section, previously starting at row 45,
# down to start at row 47).
$null = $ws.Rows.Item(46).Insert()

# Row 45 (formerly row 44, "Case File – anyone can subscribe") carries
# the formatting template we want for the new CASE_FILE row.
$null = $ws.Range("A45:G45").Copy()
$null = $ws.Range("A46:G46").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("B46").Value = "Case File – anyone can add tag"
$ws.Range("C46").Value = "CASE_FILE"
$ws.Range("G46").Value = "grant addTag to *"
$ws.Rows.Item(46).RowHeight = 30

# Match the author's final selection/viewport and persist.
$null = $ws.Range("B62").Select()
$null = $wb.Save()

Write-Host "Inserted 'grant addTag to *' rules for Complaint and Case File."
